# Book of Itza review — add a "Meta description" paragraph right after the
# title, drop the duplicate title paragraph near the end, and swap the
# trailing italic paragraph's text for the image-generation prompt.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new paragraph right after the Heading1 title with:
#       "Meta description" (bold) + ": Discover a balanced ... free."
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$boldText = "Meta description"
$restText = ": Discover a balanced and enjoyable gaming experience with innovative gameplay mechanics and lucrative bonus features in Book of Itza online slot game. Play for free."

$metaPara.Range.Text = $boldText + $restText

$metaParaAgain = $d.Paragraphs.Item(2)
$paraStart = $metaParaAgain.Range.Start
$boldRange = $d.Range($paraStart, $paraStart + $boldText.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Remove the duplicate title paragraph ("Play Book of Itza Free - ...")
#    that used to sit right before the closing italic paragraph.
# ---------------------------------------------------------------------
$old = "Play Book of Itza Free - Innovative Aztec Themed Slot"
$found = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $old -and $i -ne 1) {
        $para.Range.Delete()
        $found = $true
        break
    }
}

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-generation prompt, keeping the italic run formatting intact.
# ---------------------------------------------------------------------
$newPrompt = 'Prompt: Create a feature image for "Book of Itza" that showcases the happy Maya warrior with glasses in a cartoon style. The image should feature the warrior standing in front of an ancient temple, holding the titular Book of Itza. The temple should have Aztec designs and symbols, and the background should be vibrant and colorful. The warrior should have a big smile on his face and eye-catching details such as feathered headdress and intricate tattoos. The image should convey the excitement and adventure of playing the slot game while highlighting the Aztec theme and the expanding wild feature represented by the Book of Itza.'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$start = $lastPara.Range.Start
$end = $lastPara.Range.End - 1
$textRange = $d.Range($start, $end)
$textRange.Text = $newPrompt

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count); duplicate title removed: $found; last para now: $($d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)"
